$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "maternidade"
$ws.Range("B6").Value = -0.001692763436309752
$ws.Range("C6").Value = 0.04064422874685059
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""

$ws.Range("A7").Value = "linguagem e representação"
$ws.Range("B7").Value = 0.02099224995423199
$ws.Range("C7").Value = 0.0184696832214335
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""

$ws.Range("A8").Value = "família"
$ws.Range("B8").Value = -0.001269035532994955
$ws.Range("C8").Value = 0.0527149194004024
$ws.Range("D8").Value = -0.03103448275862064
$ws.Range("E8").Value = ""

$ws.Range("A9").Value = "resistência e luta"
$ws.Range("B9").Value = 0.01725512339323088
$ws.Range("C9").Value = 0.05009995677544854
$ws.Range("D9").Value = 0.6644219977553312
$ws.Range("E9").Value = -0.01700680272108857

$ws.Range("A10").Value = "saudade, luto ou perda"
$ws.Range("B10").Value = 0.1536106750392465
$ws.Range("C10").Value = 0.05801075583007054
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -0.01013513513513509

$ws.Range("A11").Value = "sonho e fantasia"
$ws.Range("B11").Value = 0.1261784460766218
$ws.Range("C11").Value = 0.08651128822277976
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = -0.003355704697986628

$ws.Range("A12").Value = "questão agrária e territorial"
$ws.Range("B12").Value = -0.002116850127011016
$ws.Range("C12").Value = 0.1291348610949269
$ws.Range("D12").Value = 0.7440068493150684
$ws.Range("E12").Value = -0.003355704697986628

$ws.Range("A13").Value = "mobilidade"
$ws.Range("B13").Value = 0.03837875683945757
$ws.Range("C13").Value = 0.04265942136199175
$ws.Range("D13").Value = 0.5626523647001462
$ws.Range("E13").Value = -0.003355704697986628

$ws.Range("A14").Value = "tecnologia, inovação e sociedade"
$ws.Range("B14").Value = -0.001692763436309752
$ws.Range("C14").Value = 0.1296105179647159
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = 0

$ws.Range("A15").Value = "sistema penitenciário"
$ws.Range("B15").Value = -0.0004226542688081647
$ws.Range("C15").Value = 0.2139456779045861
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0

$ws.Range("A16").Value = "vida rural, vida no campo"
$ws.Range("B16").Value = 0.1127864095278399
$ws.Range("C16").Value = 0.09277065458771383
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0

$ws.Range("A17").Value = "recreação, lazer e entretenimento"
$ws.Range("B17").Value = -0.003816793893129722
$ws.Range("C17").Value = 0.0199479651302199
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0

$ws.Range("A18").Value = "estruturas sociais e econômicas"
$ws.Range("B18").Value = 0.06184225315452274
$ws.Range("C18").Value = 0.06810972814623684
$ws.Range("D18").Value = -0.006734006734006703
$ws.Range("E18").Value = 0

$ws.Range("A19").Value = "reflexão"
$ws.Range("B19").Value = -0.003391267486222871
$ws.Range("C19").Value = 0.00322664824118124
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0.1535073922617176

$ws.Range("A20").Value = "vida cotidiana"
$ws.Range("B20").Value = 0.07433873720136508
$ws.Range("C20").Value = 0.06762693418380306
$ws.Range("D20").Value = -0.003355704697986628
$ws.Range("E20").Value = 0.200840015273005

$ws.Range("A21").Value = "cultural"
$ws.Range("B21").Value = 0.08051205906189551
$ws.Range("C21").Value = 0.07526982506934388
$ws.Range("D21").Value = 0.5993031358885017
$ws.Range("E21").Value = 0.3918644067796611

$ws.Range("A22").Value = "pandemia"
$ws.Range("B22").Value = 0.1337371292667811
$ws.Range("C22").Value = 0.1362333666309427
$ws.Range("D22").Value = 0.7367957746478873
$ws.Range("E22").Value = 0.4291714394807178

$ws.Range("A23").Value = "arte"
$ws.Range("B23").Value = 0.09271720613287915
$ws.Range("C23").Value = 0.03120625109756925
$ws.Range("D23").Value = 0.3416168606480026
$ws.Range("E23").Value = 0.4837545126353791

$ws.Range("A24").Value = "memória e patrimônio"
$ws.Range("B24").Value = 0.13437396333348
$ws.Range("C24").Value = 0.09129259470778595
$ws.Range("D24").Value = 0.200840015273005
$ws.Range("E24").Value = 0.488013698630137

$ws.Range("A25").Value = "solidão"
$ws.Range("B25").Value = 0.09216980560770316
$ws.Range("C25").Value = 0.0399854446331066
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = 0.488013698630137

$ws.Range("A26").Value = "dinâmica urbana"
$ws.Range("B26").Value = 0.1179422133191038
$ws.Range("C26").Value = 0.05848410860569164
$ws.Range("D26").Value = 0.3078703703703703
$ws.Range("E26").Value = 0.488013698630137

$ws.Range("A27").Value = "sonoridade e paisagem sonora"
$ws.Range("B27").Value = 0.07776048456563267
$ws.Range("C27").Value = 0.1093142377884295
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0.4949324324324325

$ws.Range("A28").Value = "territorialidade e colonialismo"
$ws.Range("B28").Value = 0.1686322188449849
$ws.Range("C28").Value = 0.04430741982339381
$ws.Range("D28").Value = 0.8858342878961436
$ws.Range("E28").Value = 0.4949324324324325

$ws.Range("A29").Value = "política"
$ws.Range("B29").Value = 0.0535441657579061
$ws.Range("C29").Value = 0.1991659338461489
$ws.Range("D29").Value = 0.4949324324324325
$ws.Range("E29").Value = 0.4949324324324325

$ws.Range("A30").Value = "saúde mental"
$ws.Range("B30").Value = 0.1772537902132381
$ws.Range("C30").Value = 0.07338090055341351
$ws.Range("D30").Value = 0.6088947024198823
$ws.Range("E30").Value = 0.551948051948052

$ws.Range("A31").Value = "poesia e ensaio"
$ws.Range("B31").Value = 0.0791303007391857
$ws.Range("C31").Value = 0.05965531139032709
$ws.Range("D31").Value = 0.3915343915343915
$ws.Range("E31").Value = 0.572857142857143

$ws.Range("A32").Value = "povos originários e comunidades tradicionais"
$ws.Range("B32").Value = 0.07360328350894407
$ws.Range("C32").Value = 0.1745641684215417
$ws.Range("D32").Value = 0.8542174549000487
$ws.Range("E32").Value = 0.5993031358885017

$ws.Range("A33").Value = "mulher, feminino e feminismo"
$ws.Range("B33").Value = 0.1561350573463302
$ws.Range("C33").Value = 0.1130381636868613
$ws.Range("D33").Value = 0.5297263290342875
$ws.Range("E33").Value = 0.6502923976608187

$ws.Range("A34").Value = "corpo, performance e expressão"
$ws.Range("B34").Value = 0.04792576960587136
$ws.Range("C34").Value = 0.1525330728009373
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = 0.6575028636884307

$ws.Range("A35").Value = "distopia, ficção científica e futuros imaginados"
$ws.Range("B35").Value = 0.1307555460293643
$ws.Range("C35").Value = 0.07923305428031335
$ws.Range("D35").Value = 0.6644219977553312
$ws.Range("E35").Value = 0.6609977324263039

$ws.Range("A36").Value = "violências e preconceitos de gênero"
$ws.Range("B36").Value = 0.0702501812907903
$ws.Range("C36").Value = 0.2242655108066456
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 0.6644219977553312

$ws.Range("A37").Value = "sociedade e consumo"
$ws.Range("B37").Value = 0.05153299738935113
$ws.Range("C37").Value = 0.09007331853838785
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0.6644219977553312

$ws.Range("A38").Value = "moradia e habitação"
$ws.Range("B38").Value = 0.09332284988508532
$ws.Range("C38").Value = 0.03655807124504251
$ws.Range("D38").Value = 0.3078703703703703
$ws.Range("E38").Value = 0.6644219977553312

$ws.Range("A39").Value = "violência"
$ws.Range("B39").Value = -0.003816793893129722
$ws.Range("C39").Value = 0.02542126680808454
$ws.Range("D39").Value = 0.7972881355932203
$ws.Range("E39").Value = 0.6644219977553312

$ws.Range("A40").Value = "trabalho e ofício"
$ws.Range("B40").Value = 0.1432888540031397
$ws.Range("C40").Value = 0.1128546518675327
$ws.Range("D40").Value = 0.6575028636884307
$ws.Range("E40").Value = 0.7440068493150684

$ws.Range("A41").Value = "vida afetiva"
$ws.Range("B41").Value = 0.08348697181222353
$ws.Range("C41").Value = 0.06196447104692349
$ws.Range("D41").Value = 0.3416168606480026
$ws.Range("E41").Value = 0.7440068493150684

$ws.Range("A42").Value = "alimentação e tratamentos tradicionais"
$ws.Range("B42").Value = 0.1344466114090688
$ws.Range("C42").Value = 0.230742353695167
$ws.Range("D42").Value = 0.3918644067796611
$ws.Range("E42").Value = 0.7440068493150684

$ws.Range("A43").Value = "educação e socialização"
$ws.Range("B43").Value = 0.1114948245605802
$ws.Range("C43").Value = 0.04290315008460299
$ws.Range("D43").Value = 0.5875862068965517
$ws.Range("E43").Value = 0.7440068493150684

$ws.Range("A44").Value = "desinformação, populismo e polarização"
$ws.Range("B44").Value = 0.17932651580411
$ws.Range("C44").Value = 0.2203903889728095
$ws.Range("D44").Value = -0.003355704697986628
$ws.Range("E44").Value = 0.7972881355932203

$ws.Range("A45").Value = "ambiental"
$ws.Range("B45").Value = 0.1704551937731774
$ws.Range("C45").Value = 0.1072486287876876
$ws.Range("D45").Value = 0.5814989184374602
$ws.Range("E45").Value = 0.8044473512099412

$ws.Range("A46").Value = "religião, espiritualidade e cosmologias"
$ws.Range("B46").Value = 0.1733327259475219
$ws.Range("C46").Value = 0.06693363011004605
$ws.Range("D46").Value = 0.7972881355932203
$ws.Range("E46").Value = 0.8858342878961436

$ws.Range("A47").Value = "crises e desastres ambientais e sociais"
$ws.Range("B47").Value = 0.1116421047364077
$ws.Range("C47").Value = 0.08601164557975172
$ws.Range("D47").Value = 0.2320205479452055
$ws.Range("E47").Value = 0.9059452658068575

$ws.Range("A48").Value = "pessoas com deficiência"
$ws.Range("B48").Value = 0.1761349878934625
$ws.Range("C48").Value = 0.1009122474323939
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = 1

$ws.Range("A49").Value = "gênero e sexualidade"
$ws.Range("B49").Value = 0.1065885559691625
$ws.Range("C49").Value = 0.2282058046313643
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 1

$ws.Range("A50").Value = "direitos humanos"
$ws.Range("B50").Value = 0.1381235586843064
$ws.Range("C50").Value = 0.0293775489936936
$ws.Range("D50").Value = ""
$ws.Range("E50").Value = 1

$ws.Range("A51").Value = "biografia"
$ws.Range("B51").Value = 0.0535441657579061
$ws.Range("C51").Value = 0.009196823306502511
$ws.Range("D51").Value = ""
$ws.Range("E51").Value = 1

$ws.Range("A52").Value = "amizade"
$ws.Range("B52").Value = 0.1417693981145758
$ws.Range("C52").Value = 0.2244857069255272
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 1
